{"js": "// Rewrite the body of the \"model-testing\" view template:\n//  - Heading becomes a literal \"Example Report\" title (was \"<%=Model. Title %>\").\n//  - The three ad-hoc \"First you create a view...\" / \"Things need to be\n//    inserted...\" demo paragraphs are removed.\n//  - The \"This many projects:\" / image / foreach(project) / project title /\n//    description / organizations loop block is added, followed by the\n//    closing \"<% } %>\" markers, a trailing \"test\" paragraph and a blank\n//    paragraph before the section break.\n//  - The lone \"_GoBack\" bookmark paragraph is kept (now right after the\n//    image placeholder paragraph).\n//\n// The whole body is rebuilt in one shot via Range.insertOoxml (Flat OPC),\n// which is the Office.js analogue of pasting/inserting a WordOpenXML\n// fragment \u2014 this lets us control run/paragraph boundaries, the Heading1\n// style re-use, the <w:proofErr> spell-check markers Word leaves around\n// camelCase code tokens, and the bookmark precisely.\n\nconst body = context.document.body;\n\nconst newBodyFragment =\n  '<w:p><w:pPr><w:pStyle w:val=\"Heading1\"/></w:pPr><w:r><w:t>Example Report</w:t></w:r></w:p>' +\n  '<w:p/>' +\n  '<w:p><w:r><w:t xml:space=\"preserve\">This many projects: </w:t></w:r><w:r><w:t xml:space=\"preserve\">&lt;%= </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Model.Projects.Count</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> %&gt;</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>&lt;%= Image(\\u201Ctesting.png\\u201D); %&gt;</w:t></w:r></w:p>' +\n  '<w:p><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p>' +\n  '<w:p><w:r><w:t xml:space=\"preserve\">&lt;% foreach(var project in </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Model.Projects</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>) { %&gt;</w:t></w:r></w:p>' +\n  '<w:p/>' +\n  '<w:p><w:pPr><w:pStyle w:val=\"Heading1\"/></w:pPr><w:r><w:t xml:space=\"preserve\">&lt;%= </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>project.ProjectName</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> %&gt;</w:t></w:r></w:p>' +\n  '<w:p/>' +\n  '<w:p><w:r><w:t>Description</w:t></w:r><w:r><w:t xml:space=\"preserve\">: &lt;%= </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>project.Project</w:t></w:r><w:r><w:t>Description</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> %&gt;</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>Organizations</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t xml:space=\"preserve\">&lt;% foreach(var organization in </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>project.</w:t></w:r><w:r><w:t>Organizations</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>) {</w:t></w:r><w:r><w:t>!</w:t></w:r><w:r><w:t xml:space=\"preserve\"> %&gt;</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>&lt;% } %&gt;</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>&lt;% } %&gt;</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>test</w:t></w:r></w:p>' +\n  // Word always folds the very last paragraph mark of a whole-body Replace\n  // into the body's pre-existing (unremovable) trailing paragraph, so one\n  // extra trailing empty paragraph is supplied to end up with the wanted\n  // final blank paragraph right before the section break.\n  '<w:p/>' +\n  '<w:p/>';\n\nconst flatOpcPackage =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' + newBodyFragment + '</w:body>' +\n  '</w:document>' +\n  '</pkg:xmlData>' +\n  '</pkg:part>' +\n  '</pkg:package>';\n\nbody.insertOoxml(flatOpcPackage, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Rewrite the body of the \"model-testing\" view template:\n#  - Heading becomes a literal \"Example Report\" title (was \"<%=Model. Title %>\").\n#  - The three ad-hoc \"First you create a view...\" / \"Things need to be\n#    inserted...\" demo paragraphs are removed.\n#  - The \"This many projects:\" / image / foreach(project) / project title /\n#    description / organizations loop block is added, followed by the\n#    closing \"<% } %>\" markers, a trailing \"test\" paragraph and a blank\n#    paragraph before the section break.\n#  - The lone \"_GoBack\" bookmark paragraph is kept (now right after the\n#    image placeholder paragraph).\n#\n# The whole body is rebuilt in one shot via Range.InsertXML, the COM\n# analogue of Office.js's Range.insertOoxml \u2014 this lets us control\n# run/paragraph boundaries, the Heading1 style re-use, the <w:proofErr>\n# spell-check markers Word leaves around camelCase code tokens, and the\n# bookmark precisely.\n\n$d = $word.ActiveDocument\n\n$newBodyFragment = (\n  '<w:p><w:pPr><w:pStyle w:val=\"Heading1\"/></w:pPr><w:r><w:t>Example Report</w:t></w:r></w:p>' +\n  '<w:p/>' +\n  '<w:p><w:r><w:t xml:space=\"preserve\">This many projects: </w:t></w:r><w:r><w:t xml:space=\"preserve\">&lt;%= </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Model.Projects.Count</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> %&gt;</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>&lt;%= Image(&#8220;testing.png&#8221;); %&gt;</w:t></w:r></w:p>' +\n  '<w:p><w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/><w:bookmarkEnd w:id=\"0\"/></w:p>' +\n  '<w:p><w:r><w:t xml:space=\"preserve\">&lt;% foreach(var project in </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>Model.Projects</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>) { %&gt;</w:t></w:r></w:p>' +\n  '<w:p/>' +\n  '<w:p><w:pPr><w:pStyle w:val=\"Heading1\"/></w:pPr><w:r><w:t xml:space=\"preserve\">&lt;%= </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>project.ProjectName</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> %&gt;</w:t></w:r></w:p>' +\n  '<w:p/>' +\n  '<w:p><w:r><w:t>Description</w:t></w:r><w:r><w:t xml:space=\"preserve\">: &lt;%= </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>project.Project</w:t></w:r><w:r><w:t>Description</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t xml:space=\"preserve\"> %&gt;</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>Organizations</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t xml:space=\"preserve\">&lt;% foreach(var organization in </w:t></w:r><w:proofErr w:type=\"spellStart\"/><w:r><w:t>project.</w:t></w:r><w:r><w:t>Organizations</w:t></w:r><w:proofErr w:type=\"spellEnd\"/><w:r><w:t>) {</w:t></w:r><w:r><w:t>!</w:t></w:r><w:r><w:t xml:space=\"preserve\"> %&gt;</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>&lt;% } %&gt;</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>&lt;% } %&gt;</w:t></w:r></w:p>' +\n  '<w:p><w:r><w:t>test</w:t></w:r></w:p>' +\n  '<w:p/>'\n)\n\n$xmlPackage = (\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' + $newBodyFragment + '</w:body>' +\n  '</w:document>'\n)\n\n[void]$d.Content.InsertXML($xmlPackage)\n"}
